$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2, shifting existing data down by one
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()

# Populate the newly inserted row 2 with new data
$ws.Cells.Item(2, 1).Value = -0.06342706156940325
$ws.Cells.Item(2, 2).Value = 0.2977593003249712
$ws.Cells.Item(2, 3).Value = 0.5447398664020925

# Append new rows at the bottom (rows 23-31) with new data
$ws.Cells.Item(23, 1).Value = 2.753937654378926
$ws.Cells.Item(23, 2).Value = -9.407423193861781
$ws.Cells.Item(23, 3).Value = -3.08937735819234
$ws.Cells.Item(24, 1).Value = -4.065274791019721
$ws.Cells.Item(24, 2).Value = -1.80112353185328
$ws.Cells.Item(24, 3).Value = 2.134828872796967
$ws.Cells.Item(25, 1).Value = -7.774596919373714
$ws.Cells.Item(25, 2).Value = 2.037927262666848
$ws.Cells.Item(25, 3).Value = 2.593074496199395
$ws.Cells.Item(26, 1).Value = 0.1093276535592701
$ws.Cells.Item(26, 2).Value = 10.26185343905173
$ws.Cells.Item(26, 3).Value = -3.281300154764442
$ws.Cells.Item(27, 1).Value = 3.391102220953991
$ws.Cells.Item(27, 2).Value = -5.063158106513092
$ws.Cells.Item(27, 3).Value = -1.229867340406294
$ws.Cells.Item(28, 1).Value = 6.242717754550078
$ws.Cells.Item(28, 2).Value = -0.2678701499613687
$ws.Cells.Item(28, 3).Value = -4.157948156682439
$ws.Cells.Item(29, 1).Value = 1.95433324720798
$ws.Cells.Item(29, 2).Value = -6.421389656095905
$ws.Cells.Item(29, 3).Value = 0.3005734565781801
$ws.Cells.Item(30, 1).Value = -2.448658175584859
$ws.Cells.Item(30, 2).Value = -1.071053583447552
$ws.Cells.Item(30, 3).Value = 4.803342400527646
$ws.Cells.Item(31, 1).Value = -7.789634487977738
$ws.Cells.Item(31, 2).Value = -0.2244347770038115
$ws.Cells.Item(31, 3).Value = 7.027555852401433
